# The deck currently uses the "Integral" theme (colour scheme) for its
# slide master / slides, while the (otherwise unused, notes-master-only)
# theme part still holds the default "Office Theme" colours.
#
# The authored change swaps the two colour sets: the presentation's
# active theme reverts to the stock Office theme colours (the values
# that used to live in the other theme part), i.e. the Design/Theme
# colour scheme used by the slides is changed from "Integral" back to
# the default "Office" palette.
#
# PowerPoint's COM colour properties are plain OLE_COLOR values, i.e.
# 0x00BBGGRR (blue/green/red byte order) rather than 0x00RRGGBB, so we
# build values with a small RGB() helper exactly like VBA's RGB()
# macro does.
function New-OleColor([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 0x100) + ($b * 0x10000)
}

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colors = $master.Theme.ThemeColorScheme

# MsoThemeColorSchemeIndex order: 1 dk1, 2 lt1, 3 dk2, 4 lt2,
# 5-10 accent1..accent6, 11 hlink, 12 folHlink.
$colors.Colors(1).RGB  = New-OleColor 0x00 0x00 0x00   # dk1      (unchanged)
$colors.Colors(2).RGB  = New-OleColor 0xFF 0xFF 0xFF   # lt1      (unchanged)
$colors.Colors(3).RGB  = New-OleColor 0x44 0x54 0x6A   # dk2      Integral 455F51 -> Office 44546A
$colors.Colors(4).RGB  = New-OleColor 0xE7 0xE6 0xE6   # lt2      Integral E3DED1 -> Office E7E6E6
$colors.Colors(5).RGB  = New-OleColor 0x5B 0x9B 0xD5   # accent1  Integral 99CB38 -> Office 5B9BD5
$colors.Colors(6).RGB  = New-OleColor 0xED 0x7D 0x31   # accent2  Integral 63A537 -> Office ED7D31
$colors.Colors(7).RGB  = New-OleColor 0xA5 0xA5 0xA5   # accent3  Integral E6D024 -> Office A5A5A5
$colors.Colors(8).RGB  = New-OleColor 0xFF 0xC0 0x00   # accent4  Integral CC9700 -> Office FFC000
$colors.Colors(9).RGB  = New-OleColor 0x44 0x72 0xC4   # accent5  Integral 4EB3CF -> Office 4472C4
$colors.Colors(10).RGB = New-OleColor 0x70 0xAD 0x47   # accent6  Integral 378DA6 -> Office 70AD47
$colors.Colors(11).RGB = New-OleColor 0x05 0x63 0xC1   # hlink    Integral 6B9F25 -> Office 0563C1
$colors.Colors(12).RGB = New-OleColor 0x95 0x4F 0x72   # folHlink Integral B26B02 -> Office 954F72
